# Final commit of upload excel file
# - Update a handful of contact-detail text values (FirstName, Street, Hobbies)
# - Bump the row height of the 3 header/data rows slightly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (first contact): FirstName tintu -> rohan, Hobbies trailing comma removed
$ws.Range("B2").Value = "rohan"
$ws.Range("L2").Value = "Reading ,Drawing"

# Row 3 (second contact): FirstName Maya -> mini, Street dfbdf -> abcd, Hobbies trailing comma removed
$ws.Range("B3").Value = "mini"
$ws.Range("H3").Value = "abcd"
$ws.Range("L3").Value = "Reading ,Writing"

# Slightly taller rows for the header + two data rows
$ws.Rows("1:3").RowHeight = 19.5
